# Update "想去人数" (want-to-go count) figures in column F
# for sheets "展览" and "全部类型" to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$targetValues = @{
    2  = 1104
    3  = 418
    5  = 8781
    8  = 652
    9  = 293
    10 = 159
    11 = 21
    13 = 3629
    15 = 368
    16 = 83
    17 = 2290
    20 = 315
    21 = 210
    22 = 2422
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $targetValues.Keys) {
        $ws.Range("F$row").Value = $targetValues[$row]
    }
}
